{"js": "// Replace the date heading and the 20x5 grid of arithmetic expressions\n// with their updated values. The document body flattens to 101\n// paragraphs in reading order: paragraph 0 is the date heading, and\n// paragraphs 1..100 are the table cells (row-major, 20 rows x 5 cols).\n// `oldValues`/`newValues` below are aligned to that same order.\n\nconst oldValues = [\n  \"2022-12-01 Thursday\", \"88-64=\", \"30+3=\", \"37+61=\", \"47+52=\", \"45-33=\",\n  \"70-49=\", \"41+29=\", \"22+49=\", \"64-63=\", \"80+2=\", \"98-37=\",\n  \"6+73=\", \"41+47=\", \"87-78=\", \"6+91=\", \"20-10=\", \"75-61=\",\n  \"75+3=\", \"56-28=\", \"55-10=\", \"19+46=\", \"97-49=\", \"27+2=\",\n  \"31-27=\", \"83-30=\", \"38+51=\", \"77-68=\", \"4+47=\", \"11+67=\",\n  \"15+64=\", \"51-3=\", \"94-91=\", \"82+10=\", \"19+48=\", \"19+19=\",\n  \"58-33=\", \"40-6=\", \"95-9=\", \"15+40=\", \"59-32=\", \"88-77=\",\n  \"60-43=\", \"15-14=\", \"48-42=\", \"67-16=\", \"49+30=\", \"7+34=\",\n  \"92-8=\", \"62-46=\", \"39+59=\", \"61+25=\", \"43+52=\", \"58+28=\",\n  \"43+46=\", \"80-45=\", \"6+74=\", \"46-27=\", \"34+27=\", \"97-93=\",\n  \"4+26=\", \"21+55=\", \"34+28=\", \"50-7=\", \"43+4=\", \"91-58=\",\n  \"65-44=\", \"41+29=\", \"78+17=\", \"63-60=\", \"68-50=\", \"73-15=\",\n  \"16+57=\", \"0+52=\", \"3+2=\", \"39+50=\", \"25-9=\", \"52-3=\",\n  \"92-9=\", \"55-35=\", \"18+49=\", \"5+28=\", \"78-29=\", \"20-13=\",\n  \"37+9=\", \"84+10=\", \"74-68=\", \"61-9=\", \"14+0=\", \"58+13=\",\n  \"98-43=\", \"27+5=\", \"64-10=\", \"40+6=\", \"98-10=\", \"54-50=\",\n  \"18+49=\", \"3+59=\", \"68-39=\", \"58-32=\", \"3+67=\"\n];\n\nconst newValues = [\n  \"2022-12-02 Friday\", \"85-79=\", \"65-53=\", \"34+57=\", \"16+54=\", \"43+43=\",\n  \"38-38=\", \"89+1=\", \"36-4=\", \"80-9=\", \"64-39=\", \"61-28=\",\n  \"86-14=\", \"25+44=\", \"58-16=\", \"82-48=\", \"88-58=\", \"73-60=\",\n  \"10+43=\", \"56+11=\", \"44+45=\", \"99-0=\", \"52+11=\", \"45-44=\",\n  \"76-31=\", \"39+18=\", \"15+63=\", \"82-12=\", \"6+39=\", \"19+56=\",\n  \"34+23=\", \"49-49=\", \"2+22=\", \"78-4=\", \"82-19=\", \"92-31=\",\n  \"73-22=\", \"61+17=\", \"9-0=\", \"37+14=\", \"50+19=\", \"79-20=\",\n  \"5+23=\", \"54-37=\", \"83-44=\", \"64+32=\", \"62-13=\", \"66-42=\",\n  \"66-14=\", \"66-26=\", \"62-45=\", \"32+63=\", \"99-56=\", \"17+50=\",\n  \"66-3=\", \"64+2=\", \"69+23=\", \"12+39=\", \"11+39=\", \"10+20=\",\n  \"7+79=\", \"53+8=\", \"90-80=\", \"74-66=\", \"30-23=\", \"41-24=\",\n  \"36-2=\", \"38-8=\", \"14+25=\", \"33+7=\", \"88-67=\", \"26-3=\",\n  \"28+43=\", \"56+24=\", \"58-12=\", \"90-32=\", \"77+3=\", \"11+5=\",\n  \"53+44=\", \"10+48=\", \"13+60=\", \"78-17=\", \"69+2=\", \"80-56=\",\n  \"95-26=\", \"81-63=\", \"48-47=\", \"16+59=\", \"69+15=\", \"12-2=\",\n  \"98-63=\", \"52+46=\", \"9+48=\", \"72-35=\", \"78-3=\", \"66-1=\",\n  \"0+88=\", \"88-41=\", \"38+48=\", \"38-35=\", \"7+58=\"\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nif (items.length !== oldValues.length) {\n  throw new Error(\n    \"Unexpected paragraph count: \" + items.length + \" vs expected \" + oldValues.length\n  );\n}\n\nfor (let i = 0; i < items.length; i++) {\n  const current = items[i].text;\n  if (current === oldValues[i] && current !== newValues[i]) {\n    items[i].insertText(newValues[i], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and the 20x5 grid of arithmetic expressions.\n# Word COM Range.Text for a paragraph/cell includes trailing control\n# characters (paragraph mark \\r / cell mark \\a), so trim those before\n# comparing against the plain text we expect to see.\n\n$d = $word.ActiveDocument\n\nfunction Get-PlainText($range) {\n    return $range.Text.TrimEnd([char]13, [char]7)\n}\n\n# --- Update the date heading paragraph ---\n$dateOld = \"2022-12-01 Thursday\"\n$dateNew = \"2022-12-02 Friday\"\n$p1 = $d.Paragraphs.Item(1)\n$p1Text = Get-PlainText $p1.Range\nif ($p1Text -eq $dateOld -and $p1Text -ne $dateNew) {\n    $p1.Range.Text = $dateNew\n}\n\n# --- Update the table of arithmetic expressions (row-major order) ---\n$oldValues = @(\n    \"88-64=\",\n    \"30+3=\",\n    \"37+61=\",\n    \"47+52=\",\n    \"45-33=\",\n    \"70-49=\",\n    \"41+29=\",\n    \"22+49=\",\n    \"64-63=\",\n    \"80+2=\",\n    \"98-37=\",\n    \"6+73=\",\n    \"41+47=\",\n    \"87-78=\",\n    \"6+91=\",\n    \"20-10=\",\n    \"75-61=\",\n    \"75+3=\",\n    \"56-28=\",\n    \"55-10=\",\n    \"19+46=\",\n    \"97-49=\",\n    \"27+2=\",\n    \"31-27=\",\n    \"83-30=\",\n    \"38+51=\",\n    \"77-68=\",\n    \"4+47=\",\n    \"11+67=\",\n    \"15+64=\",\n    \"51-3=\",\n    \"94-91=\",\n    \"82+10=\",\n    \"19+48=\",\n    \"19+19=\",\n    \"58-33=\",\n    \"40-6=\",\n    \"95-9=\",\n    \"15+40=\",\n    \"59-32=\",\n    \"88-77=\",\n    \"60-43=\",\n    \"15-14=\",\n    \"48-42=\",\n    \"67-16=\",\n    \"49+30=\",\n    \"7+34=\",\n    \"92-8=\",\n    \"62-46=\",\n    \"39+59=\",\n    \"61+25=\",\n    \"43+52=\",\n    \"58+28=\",\n    \"43+46=\",\n    \"80-45=\",\n    \"6+74=\",\n    \"46-27=\",\n    \"34+27=\",\n    \"97-93=\",\n    \"4+26=\",\n    \"21+55=\",\n    \"34+28=\",\n    \"50-7=\",\n    \"43+4=\",\n    \"91-58=\",\n    \"65-44=\",\n    \"41+29=\",\n    \"78+17=\",\n    \"63-60=\",\n    \"68-50=\",\n    \"73-15=\",\n    \"16+57=\",\n    \"0+52=\",\n    \"3+2=\",\n    \"39+50=\",\n    \"25-9=\",\n    \"52-3=\",\n    \"92-9=\",\n    \"55-35=\",\n    \"18+49=\",\n    \"5+28=\",\n    \"78-29=\",\n    \"20-13=\",\n    \"37+9=\",\n    \"84+10=\",\n    \"74-68=\",\n    \"61-9=\",\n    \"14+0=\",\n    \"58+13=\",\n    \"98-43=\",\n    \"27+5=\",\n    \"64-10=\",\n    \"40+6=\",\n    \"98-10=\",\n    \"54-50=\",\n    \"18+49=\",\n    \"3+59=\",\n    \"68-39=\",\n    \"58-32=\",\n    \"3+67=\"\n)\n$newValues = @(\n    \"85-79=\",\n    \"65-53=\",\n    \"34+57=\",\n    \"16+54=\",\n    \"43+43=\",\n    \"38-38=\",\n    \"89+1=\",\n    \"36-4=\",\n    \"80-9=\",\n    \"64-39=\",\n    \"61-28=\",\n    \"86-14=\",\n    \"25+44=\",\n    \"58-16=\",\n    \"82-48=\",\n    \"88-58=\",\n    \"73-60=\",\n    \"10+43=\",\n    \"56+11=\",\n    \"44+45=\",\n    \"99-0=\",\n    \"52+11=\",\n    \"45-44=\",\n    \"76-31=\",\n    \"39+18=\",\n    \"15+63=\",\n    \"82-12=\",\n    \"6+39=\",\n    \"19+56=\",\n    \"34+23=\",\n    \"49-49=\",\n    \"2+22=\",\n    \"78-4=\",\n    \"82-19=\",\n    \"92-31=\",\n    \"73-22=\",\n    \"61+17=\",\n    \"9-0=\",\n    \"37+14=\",\n    \"50+19=\",\n    \"79-20=\",\n    \"5+23=\",\n    \"54-37=\",\n    \"83-44=\",\n    \"64+32=\",\n    \"62-13=\",\n    \"66-42=\",\n    \"66-14=\",\n    \"66-26=\",\n    \"62-45=\",\n    \"32+63=\",\n    \"99-56=\",\n    \"17+50=\",\n    \"66-3=\",\n    \"64+2=\",\n    \"69+23=\",\n    \"12+39=\",\n    \"11+39=\",\n    \"10+20=\",\n    \"7+79=\",\n    \"53+8=\",\n    \"90-80=\",\n    \"74-66=\",\n    \"30-23=\",\n    \"41-24=\",\n    \"36-2=\",\n    \"38-8=\",\n    \"14+25=\",\n    \"33+7=\",\n    \"88-67=\",\n    \"26-3=\",\n    \"28+43=\",\n    \"56+24=\",\n    \"58-12=\",\n    \"90-32=\",\n    \"77+3=\",\n    \"11+5=\",\n    \"53+44=\",\n    \"10+48=\",\n    \"13+60=\",\n    \"78-17=\",\n    \"69+2=\",\n    \"80-56=\",\n    \"95-26=\",\n    \"81-63=\",\n    \"48-47=\",\n    \"16+59=\",\n    \"69+15=\",\n    \"12-2=\",\n    \"98-63=\",\n    \"52+46=\",\n    \"9+48=\",\n    \"72-35=\",\n    \"78-3=\",\n    \"66-1=\",\n    \"0+88=\",\n    \"88-41=\",\n    \"38+48=\",\n    \"38-35=\",\n    \"7+58=\"\n)\n\n$t = $d.Tables.Item(1)\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\n\nif (($rows * $cols) -ne $oldValues.Count) {\n    throw \"Unexpected table size: $rows x $cols vs expected $($oldValues.Count) cells\"\n}\n\n$idx = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $cell = $t.Cell($r, $c)\n        $current = Get-PlainText $cell.Range\n        $expectedOld = $oldValues[$idx]\n        $expectedNew = $newValues[$idx]\n        if ($current -eq $expectedOld -and $current -ne $expectedNew) {\n            $cell.Range.Text = $expectedNew\n        }\n        $idx = $idx + 1\n    }\n}\n"}
